# The deck ships with two theme parts:
#   ppt/theme/theme1.xml  -> bound to the (single) Slide Master / Design,
#                             i.e. the theme that actually paints every
#                             slide. Currently the "Integral" / "Red Violet"
#                             palette (accent colours E32D91, C830CC, ...).
#   ppt/theme/theme2.xml  -> bound only to the Notes Master, already the
#                             stock "Office Theme" palette.
#
# The target edit swaps the two palettes, so the slide design becomes the
# stock "Office Theme" colours (44546A / E7E6E6 / 5B9BD5 / ED7D31 / A5A5A5 /
# FFC000 / 4472C4 / 70AD47 / 0563C1 / 954F72, with dk1/lt1 unchanged).
# Re-colour the live theme through the Design's ThemeColorScheme: every
# slide's ThemeColorScheme is backed by the one Design/SlideMaster theme
# that drives the whole deck, so this repaints theme1.xml in place.

function HexColor([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$colorScheme = $s.ThemeColorScheme

# Office Theme colour scheme, in the fixed 12-slot COM order:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
# 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeTheme = @("000000", "FFFFFF", "44546A", "E7E6E6", "5B9BD5", "ED7D31", "A5A5A5", "FFC000", "4472C4", "70AD47", "0563C1", "954F72")

for ($i = 0; $i -lt $officeTheme.Count; $i++) {
    $colorScheme.Colors($i + 1).RGB = HexColor $officeTheme[$i]
}
